# Add new columns I (I0) and J (IF) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows (row -> [I, J])
$data = @{
    2  = @(5, 6)
    3  = @(8, 8)
    4  = @(6, 7)
    5  = @(8, 8)
    6  = @(9, 9)
    7  = @(7, 9)
    8  = @(5, 7)
    9  = @(6, 7)
    10 = @(5, 7)
    11 = @(8, 8)
    12 = @(10, 10)
    13 = @(8, 8)
    14 = @(8, 8)
    15 = @(8, 8)
    16 = @(9, 9)
    17 = @(12, 12)
    18 = @(8, 8)
    19 = @(9, 9)
    20 = @(9, 9)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
